# Weekly fruit/vegetable price update: a new week of "Pepino dulce" price
# data (date serial 45077) is inserted right after the existing row 609
# block, shifting all subsequent rows down by 4 (and growing the used
# range from R672 to R676).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows starting at row 610 (pushes 610:672 -> 614:676).
$ws.Rows.Item(610).Insert()
$ws.Rows.Item(610).Insert()
$ws.Rows.Item(610).Insert()
$ws.Rows.Item(610).Insert()

# Common values shared by every data row in this sheet.
$mercadoId = 2
$mercado = "Comercializadora del Agro de Limarí"
$region = "Coquimbo"
$codreg = 4
$categoriaId = 100112043
$categoria = "Pepino dulce"
$variedad = "Cultivar IV Región"
$unidad = "`$/bandeja 18 kilos"
$origen = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

function Set-PrecioRow {
    param($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $PrecioKg)

    $ws.Cells.Item($Row, 1).Value = $mercadoId
    $ws.Cells.Item($Row, 2).Value = $mercado
    $ws.Cells.Item($Row, 3).Value = $region
    $ws.Cells.Item($Row, 4).Value = $Fecha
    $ws.Cells.Item($Row, 5).Value = $codreg
    $ws.Cells.Item($Row, 6).Value = $categoriaId
    $ws.Cells.Item($Row, 7).Value = $categoria
    $ws.Cells.Item($Row, 8).Value = $variedad
    $ws.Cells.Item($Row, 9).Value = $Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin
    $ws.Cells.Item($Row, 12).Value = $PrecioMax
    $ws.Cells.Item($Row, 13).Value = $PrecioProm
    $ws.Cells.Item($Row, 14).Value = $unidad
    $ws.Cells.Item($Row, 15).Value = $origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg
    $ws.Cells.Item($Row, 17).Value = $kgUnidades
    $ws.Cells.Item($Row, 18).Value = $clasificacion
}

Set-PrecioRow 610 45077 "Especial" 350  12000 12000 12000 667
Set-PrecioRow 611 45077 "Primera"  1250 9000  11000 9920  551
Set-PrecioRow 612 45077 "Segunda"  600  7000  8000  7500  417
Set-PrecioRow 613 45077 "Tercera"  400  5000  6000  5500  306
